# Add a new "2022-Q1" sheet (holding the quarter's fund-holding detail) right
# before the existing "总计" (totals) sheet, and add a corresponding summary
# row at the top of "总计"'s data table.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q3")

# --- 1. Create the new "2022-Q1" sheet by duplicating an existing quarter
#        sheet (so sheetPr / pageMargins / default styles all match the
#        sibling sheets) and place it immediately before "总计". -----------
$templateSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q3 (2)")
$newSheet.Name = "2022-Q1"

# Re-resolve "总计" by name: inserting the new sheet shifted everyone's
# positional index, so the earlier $totalSheet reference now points at the
# wrong tab.
$totalSheet = $wb.Worksheets.Item("总计")

# Clear out the copied sample data rows (2 onward) but keep header/styles.
$newSheet.Rows.Item(2).Resize(1).ClearContents()

# This quarter has 3 holdings (rows 2-4); rows 3-4 are brand new, so give
# column A the same index-cell style used by row 2 (and every other sheet).
$newSheet.Range("A2").Copy()
$newSheet.Range("A3:A4").PasteSpecial(-4122)

# --- 2. Fill the "2022-Q1" sheet with this quarter's fund holdings. -------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B-G on the data rows hold text (codes / formatted numbers), so
# force a text format before assigning, then restore the default style
# (matches how the sibling quarter sheets store these values).
$dataTextRange = $newSheet.Range("B2:G4")
$dataTextRange.NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "002666"
$newSheet.Range("C2").Value = "前海开源沪港深创新成长灵活配置混合A"
$newSheet.Range("D2").Value = "11.96"
$newSheet.Range("E2").Value = "81.64"
$newSheet.Range("F2").Value = "6.97"
$newSheet.Range("G2").Value = "0.8336"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "163801"
$newSheet.Range("C3").Value = "中银中国混合(LOF)"
$newSheet.Range("D3").Value = "10.14"
$newSheet.Range("E3").Value = "89.19"
$newSheet.Range("F3").Value = "2.45"
$newSheet.Range("G3").Value = "0.2484"
$newSheet.Range("H3").Value = 10

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "002667"
$newSheet.Range("C4").Value = "前海开源沪港深创新成长灵活配置混合C"
$newSheet.Range("D4").Value = "3.25"
$newSheet.Range("E4").Value = "81.64"
$newSheet.Range("F4").Value = "6.97"
$newSheet.Range("G4").Value = "0.2265"
$newSheet.Range("H4").Value = 4

$dataTextRange.Style = "Normal"

# --- 3. Insert a new top data row in "总计" for 2022-Q1, pushing the rest
#        of the table down by one row (copy preserves per-row styles). ----
$totalSheet.Range("A2:D5").Copy($totalSheet.Range("A3:D6"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 1.31

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

# Restore the original active sheet/selection.
$wb.Worksheets.Item("2020-Q4").Activate()
